$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 520, shifting rows 520:626 down to 521:627
$ws.Rows.Item(520).Insert()

# Populate the new row 520 with values (copy of former row 520, with D and J updated)
$ws.Cells.Item(520, 1).Value = 4
$ws.Cells.Item(520, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(520, 3).Value = "Los Lagos"
$ws.Cells.Item(520, 4).Value = 45244
$ws.Cells.Item(520, 5).Value = 10
$ws.Cells.Item(520, 6).Value = 100112023
$ws.Cells.Item(520, 7).Value = "Brócoli"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 1500
$ws.Cells.Item(520, 11).Value = 1500
$ws.Cells.Item(520, 12).Value = 1500
$ws.Cells.Item(520, 13).Value = 1500
$ws.Cells.Item(520, 14).Value = "$/unidad"
$ws.Cells.Item(520, 15).Value = "Región Metropolitana"
$ws.Cells.Item(520, 16).Value = 1500
$ws.Cells.Item(520, 17).Value = 1
$ws.Cells.Item(520, 18).Value = "Hortaliza"
